$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.909.96'
$ws.Range('E2').Value = '  -3.49%  '
$ws.Range('D3').Value = '3.286.76'
$ws.Range('E3').Value = '  -3.97%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '553.51'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.36%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '140.58'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -8.31%  '
$ws.Range('D8').Value = '3.283.67'
$ws.Range('E8').Value = '  -4.05%  '
$ws.Range('E9').Value = '  -3.89%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.75'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.87%  '
$ws.Range('E11').Value = '  -4.85%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.406'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.81%  '
$ws.Range('D13').Value = '3.843.27'
$ws.Range('E13').Value = '  -4.02%  '
$ws.Range('E14').Value = '  -0.20%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.87'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -6.65%  '
$ws.Range('D16').Value = '3.277.44'
$ws.Range('E16').Value = '  -4.22%  '
$ws.Range('E17').Value = '  -5.01%  '
$ws.Range('D18').Value = '60.000.67'
$ws.Range('E18').Value = '  -3.39%  '
$ws.Range('E19').Value = '  -5.93%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.80'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.10%  '
$ws.Range('E21').Value = '  -4.82%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '372.01'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '73.49'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.39%  '
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.533'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -6.55%  '
$ws.Range('D26').Value = '3.421.54'
$ws.Range('E26').Value = '  -3.96%  '
$ws.Range('E27').Value = '  -10.25%  '
$ws.Range('E28').Value = '  -5.26%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.996'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.39%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.09'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -8.13%  '
$ws.Range('E31').Value = '  +0.00%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.02'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.94%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.47'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -5.55%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '22.52'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.17%  '
$ws.Range('E35').Value = '  -7.05%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.06'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -8.30%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '166.25'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.26%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.52'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.83%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.64'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.24%  '
$ws.Range('D40').Value = '3.313.95'
$ws.Range('E40').Value = '  -4.02%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '26.13'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -16.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0726'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -7.43%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '41.65'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.63%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.11'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.83%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.57'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.85%  '
$ws.Range('E47').Value = '  -6.28%  '
$ws.Range('E48').Value = '  -0.06%  '
$ws.Range('D49').Value = '2.332.14'
$ws.Range('E49').Value = '  -8.50%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.36'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -7.36%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '21.19'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -6.25%  '
